$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.600.59'
$ws.Range("E2").Value = '  -1.76%  '
$ws.Range("D3").Value = '3.325.38'
$ws.Range("E3").Value = '  -2.60%  '
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").Value = "'579.52"
$ws.Range("E5").Value = '  -2.76%  '
$ws.Range("D6").Value = "'173.70"
$ws.Range("E6").Value = '  -7.71%  '
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = "'0.586"
$ws.Range("E8").Value = '  -2.66%  '
$ws.Range("D9").Value = '3.319.48'
$ws.Range("E9").Value = '  -2.57%  '
$ws.Range("D10").Value = "'0.175"
$ws.Range("E10").Value = '  -5.75%  '
$ws.Range("D11").Value = "'0.575"
$ws.Range("E11").Value = '  -3.05%  '
$ws.Range("D12").Value = "'45.33"
$ws.Range("E12").Value = '  -5.27%  '
$ws.Range("D13").Value = "'0.0000268"
$ws.Range("E13").Value = '  -4.89%  '
$ws.Range("D14").Value = "'667.42"
$ws.Range("E14").Value = '  +3.43%  '
$ws.Range("D15").Value = '3.857.29'
$ws.Range("E15").Value = '  -2.71%  '
$ws.Range("D16").Value = "'8.38"
$ws.Range("E16").Value = '  -3.21%  '
$ws.Range("D17").Value = '67.716.72'
$ws.Range("E17").Value = '  -1.85%  '
$ws.Range("E18").Value = '  -1.34%  '
$ws.Range("D19").Value = '3.320.81'
$ws.Range("E19").Value = '  -2.91%  '
$ws.Range("D20").Value = "'17.40"
$ws.Range("E20").Value = '  -4.31%  '
$ws.Range("D21").Value = "'10.90"
$ws.Range("E21").Value = '  -2.64%  '
$ws.Range("D22").Value = "'0.887"
$ws.Range("E22").Value = '  -3.38%  '
$ws.Range("E23").Value = '  +5.06%  '
$ws.Range("D24").Value = "'16.94"
$ws.Range("E24").Value = '  -6.31%  '
$ws.Range("D25").Value = "'98.28"
$ws.Range("E25").Value = '  -2.34%  '
$ws.Range("D26").Value = "'3.83"
$ws.Range("D27").Value = "'2.67"
$ws.Range("E27").Value = '  -7.23%  '
$ws.Range("D28").Value = "'9.28"
$ws.Range("E28").Value = '  -5.77%  '
$ws.Range("D29").Value = "'33.80"
$ws.Range("E29").Value = '  +2.46%  '
$ws.Range("B30").Value = 'NEARProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D30").Value = "'7.36"
$ws.Range("E30").Value = '  +6.66%  '
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").Value = "'8.39"
$ws.Range("E31").Value = '  -4.18%  '
$ws.Range("D32").Value = "'588.18"
$ws.Range("E32").Value = '  -4.53%  '
$ws.Range("D33").Value = "'10.90"
$ws.Range("E33").Value = '  -2.60%  '
$ws.Range("D34").Value = "'0.104"
$ws.Range("E34").Value = '  -2.72%  '
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("D36").Value = '3.695.33'
$ws.Range("E36").Value = '  -8.78%  '
$ws.Range("D37").Value = "'56.75"
$ws.Range("E37").Value = '  -0.54%  '
$ws.Range("D38").Value = "'3.26"
$ws.Range("E38").Value = '  -15.65%  '
$ws.Range("D39").Value = "'0.131"
$ws.Range("E39").Value = '  -1.03%  '
$ws.Range("D40").Value = "'32.86"
$ws.Range("E40").Value = '  -3.10%  '
$ws.Range("D41").Value = "'2.62"
$ws.Range("E41").Value = '  -7.08%  '
$ws.Range("D42").Value = "'3.10"
$ws.Range("E42").Value = '  -7.01%  '
$ws.Range("D43").Value = "'0.332"
$ws.Range("E43").Value = '  -4.30%  '
$ws.Range("B44").Value = 'PEPE'
$ws.Range("C44").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D44").Value = '0.0₃0658'
$ws.Range("E44").Value = '  -8.07%  '
$ws.Range("B45").Value = 'ApeXProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D45").Value = "'3.24"
$ws.Range("E45").Value = '  -5.31%  '
$ws.Range("D46").Value = "'0.0404"
$ws.Range("E46").Value = '  -5.38%  '
$ws.Range("B47").Value = 'Stellar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D47").Value = "'0.127"
$ws.Range("E47").Value = '  -2.40%  '
$ws.Range("B48").Value = 'ThetaToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D48").Value = "'2.57"
$ws.Range("E48").Value = '  -2.00%  '
$ws.Range("E49").Value = '  -0.14%  '
$ws.Range("E50").Value = '  -4.32%  '
$ws.Range("D51").Value = "'126.78"
$ws.Range("E51").Value = '  -2.47%  '

$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D51").Style = "Normal"
